$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value2 = 5685.2666
$ws.Range("J62").Value2 = 5730.875
$ws.Range("L62").Value2 = 5730.875
$ws.Range("N62").Value2 = -6978.875
$ws.Range("H65").Value2 = 5685.2666
$ws.Range("J65").Value2 = 5730.875
$ws.Range("L65").Value2 = 28654.375
$ws.Range("N65").Value2 = -34894.375
$ws.Range("H129").Value2 = 2894
$ws.Range("I129").Value2 = 2630.75
$ws.Range("K129").Value2 = 7892.25
$ws.Range("M129").Value2 = -2892.25
$ws.Range("H138").Value2 = 18873896
$ws.Range("I138").Value2 = 58825770
$ws.Range("J138").Value2 = 7736.1665
$ws.Range("K138").Value2 = 176477310
$ws.Range("L138").Value2 = 23208.4995
$ws.Range("M138").Value2 = -176472170
$ws.Range("N138").Value2 = -33488.49950000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 841475.9
$ws.Range("I2").Value2 = 922728
$ws.Range("K2").Value2 = 922728
$ws.Range("M2").Value2 = -922615
$ws.Range("H32").Value2 = 2664.8481
$ws.Range("I32").Value2 = 1975.4166
$ws.Range("J32").Value2 = 9756.143
$ws.Range("K32").Value2 = 1975.4166
$ws.Range("L32").Value2 = 9756.143
$ws.Range("M32").Value2 = -1688.4166
$ws.Range("N32").Value2 = -10330.143
$ws.Range("H69").Value2 = 189994
$ws.Range("J69").Value2 = 189994
$ws.Range("L69").Value2 = 189994
$ws.Range("N69").Value2 = -191492
$ws.Range("H72").Value2 = 189994
$ws.Range("J72").Value2 = 189994
$ws.Range("L72").Value2 = 569982
$ws.Range("N72").Value2 = -577470
$ws.Range("H74").Value2 = 3505.818
$ws.Range("I74").Value2 = 3645.5
$ws.Range("K74").Value2 = 3645.5
$ws.Range("M74").Value2 = -2771.5
$ws.Range("H77").Value2 = 3505.818
$ws.Range("I77").Value2 = 3645.5
$ws.Range("K77").Value2 = 18227.5
$ws.Range("M77").Value2 = -13859.5
$ws.Range("H96").Value2 = 32791.855
$ws.Range("J96").Value2 = 32791.855
$ws.Range("L96").Value2 = 32791.855
$ws.Range("N96").Value2 = -38283.855
$ws.Range("H106").Value2 = 0
$ws.Range("J106").Value2 = 0
$ws.Range("L106").Value2 = 0
$ws.Range("N106").ClearContents()
$ws.Range("H116").Value2 = 841475.9
$ws.Range("I116").Value2 = 922728
$ws.Range("K116").Value2 = 922728
$ws.Range("M116").Value2 = -920434
$ws.Range("H132").Value2 = 4448.849
$ws.Range("I132").Value2 = 4385.4346
$ws.Range("K132").Value2 = 13156.3038
$ws.Range("M132").Value2 = -10626.3038

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 841475.9
$ws.Range("I3").Value2 = 922728
$ws.Range("K3").Value2 = 922728
$ws.Range("M3").Value2 = -922614
$ws.Range("H70").Value2 = 169838
$ws.Range("J70").Value2 = 169838
$ws.Range("L70").Value2 = 169838
$ws.Range("N70").Value2 = -170424
$ws.Range("H73").Value2 = 169838
$ws.Range("J73").Value2 = 169838
$ws.Range("L73").Value2 = 169838
$ws.Range("N73").Value2 = -171866
$ws.Range("H108").Value2 = 49999.91
$ws.Range("J108").Value2 = 49999.91
$ws.Range("L108").Value2 = 49999.91
$ws.Range("N108").Value2 = -57679.91

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 148.09091
$ws.Range("I7").Value2 = 62.9
$ws.Range("K7").Value2 = 62.9
$ws.Range("M7").Value2 = 50.1
$ws.Range("H28").Value2 = 5264.1665
$ws.Range("I28").Value2 = 2500
$ws.Range("K28").Value2 = 2500
$ws.Range("M28").Value2 = -2255
$ws.Range("H43").Value2 = 550000
$ws.Range("J43").Value2 = 550000
$ws.Range("L43").Value2 = 550000
$ws.Range("N43").Value2 = -550368
$ws.Range("H101").Value2 = 550000
$ws.Range("J101").Value2 = 550000
$ws.Range("L101").Value2 = 550000
$ws.Range("N101").Value2 = -556490
$ws.Range("H132").Value2 = 6791.2964
$ws.Range("I132").Value2 = 6714.6
$ws.Range("K132").Value2 = 20143.8
$ws.Range("M132").Value2 = -17613.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value2 = 1053455.4
$ws.Range("J98").Value2 = 888.7143
$ws.Range("L98").Value2 = 2666.1429
$ws.Range("N98").Value2 = -5662.1429
$ws.Range("H132").Value2 = 1619.3077
$ws.Range("I132").Value2 = 1699
$ws.Range("J132").Value2 = 1569.5
$ws.Range("K132").Value2 = 15291
$ws.Range("L132").Value2 = 14125.5
$ws.Range("M132").Value2 = -12761
$ws.Range("N132").Value2 = -19185.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value2 = 670.4737
$ws.Range("I97").Value2 = 611.1429000000001
$ws.Range("J97").Value2 = 836.6
$ws.Range("K97").Value2 = 611.1429000000001
$ws.Range("L97").Value2 = 836.6
$ws.Range("M97").Value2 = -115.1429000000001
$ws.Range("N97").Value2 = -1828.6
$ws.Range("H113").Value2 = 1887.9231
$ws.Range("I113").Value2 = 1462
$ws.Range("K113").Value2 = 1462
$ws.Range("M113").Value2 = 708
$ws.Range("H122").Value2 = 2025.6774
$ws.Range("I122").Value2 = 1870.1765
$ws.Range("K122").Value2 = 5610.529500000001
$ws.Range("M122").Value2 = -3160.529500000001
$ws.Range("H132").Value2 = 4703.8823
$ws.Range("I132").Value2 = 3297.2273
$ws.Range("K132").Value2 = 9891.6819
$ws.Range("M132").Value2 = -7361.6819

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value2 = 22347
$ws.Range("I32").Value2 = 3513
$ws.Range("K32").Value2 = 3513
$ws.Range("M32").Value2 = -3196
$ws.Range("H40").Value2 = 5848.875
$ws.Range("I40").Value2 = 2266.3333
$ws.Range("K40").Value2 = 2266.3333
$ws.Range("M40").Value2 = -2130.3333
$ws.Range("H68").Value2 = 4036.1538
$ws.Range("I68").Value2 = 2999.7144
$ws.Range("J68").Value2 = 5245.3335
$ws.Range("K68").Value2 = 2999.7144
$ws.Range("L68").Value2 = 5245.3335
$ws.Range("M68").Value2 = -2250.7144
$ws.Range("N68").Value2 = -6743.3335
$ws.Range("H71").Value2 = 4036.1538
$ws.Range("I71").Value2 = 2999.7144
$ws.Range("J71").Value2 = 5245.3335
$ws.Range("K71").Value2 = 14998.572
$ws.Range("L71").Value2 = 26226.6675
$ws.Range("M71").Value2 = -11254.572
$ws.Range("N71").Value2 = -33714.6675
$ws.Range("H82").Value2 = 599.2973
$ws.Range("I82").Value2 = 556.8219
$ws.Range("K82").Value2 = 556.8219
$ws.Range("M82").Value2 = -195.8219
$ws.Range("H85").Value2 = 599.2973
$ws.Range("I85").Value2 = 556.8219
$ws.Range("K85").Value2 = 556.8219
$ws.Range("M85").Value2 = 691.1781
$ws.Range("H122").Value2 = 7491.5
$ws.Range("I122").Value2 = 6753.5386
$ws.Range("K122").Value2 = 20260.6158
$ws.Range("M122").Value2 = -17810.6158
$ws.Range("H136").Value2 = 4291351.5
$ws.Range("I136").Value2 = 5627868
$ws.Range("K136").Value2 = 16883604
$ws.Range("M136").Value2 = -16881054

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value2 = 72906.21000000001
$ws.Range("I81").Value2 = 2449.6667
$ws.Range("K81").Value2 = 4899.3334
$ws.Range("M81").Value2 = -3838.3334
$ws.Range("H84").Value2 = 72906.21000000001
$ws.Range("I84").Value2 = 2449.6667
$ws.Range("K84").Value2 = 24496.667
$ws.Range("M84").Value2 = -19192.667
$ws.Range("H122").Value2 = 12865.412
$ws.Range("I122").Value2 = 6081.1665
$ws.Range("J122").Value2 = 29147.6
$ws.Range("K122").Value2 = 18243.4995
$ws.Range("L122").Value2 = 87442.79999999999
$ws.Range("M122").Value2 = -15793.4995
$ws.Range("N122").Value2 = -92342.79999999999
$ws.Range("H132").Value2 = 5305.729
$ws.Range("I132").Value2 = 5061.2954
$ws.Range("K132").Value2 = 15183.8862
$ws.Range("M132").Value2 = -12653.8862
